$d = $word.ActiveDocument

$d.Content.Find.Execute("90÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷8=", 2) | Out-Null
$d.Content.Find.Execute("59÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷8=", 2) | Out-Null
$d.Content.Find.Execute("63÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷8=", 2) | Out-Null
$d.Content.Find.Execute("62÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=", 2) | Out-Null
$d.Content.Find.Execute("31÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷3=", 2) | Out-Null
$d.Content.Find.Execute("63÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷4=", 2) | Out-Null
$d.Content.Find.Execute("70÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=", 2) | Out-Null
$d.Content.Find.Execute("50÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷8=", 2) | Out-Null
$d.Content.Find.Execute("60÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=", 2) | Out-Null
$d.Content.Find.Execute("66÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷8=", 2) | Out-Null
$d.Content.Find.Execute("83÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷6=", 2) | Out-Null
$d.Content.Find.Execute("86÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷3=", 2) | Out-Null
$d.Content.Find.Execute("91÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷8=", 2) | Out-Null
$d.Content.Find.Execute("61÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷5=", 2) | Out-Null
$d.Content.Find.Execute("82÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "43÷6=", 2) | Out-Null
$d.Content.Find.Execute("27÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷8=", 2) | Out-Null
$d.Content.Find.Execute("96÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷9=", 2) | Out-Null
$d.Content.Find.Execute("36÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷9=", 2) | Out-Null
$d.Content.Find.Execute("31÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷4=", 2) | Out-Null
$d.Content.Find.Execute("80÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷9=", 2) | Out-Null
$d.Content.Find.Execute("92÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷2=", 2) | Out-Null
$d.Content.Find.Execute("78÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷4=", 2) | Out-Null
$d.Content.Find.Execute("16÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷7=", 2) | Out-Null
$d.Content.Find.Execute("57÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷2=", 2) | Out-Null
$d.Content.Find.Execute("54÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷9=", 2) | Out-Null
